$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (E1 -> "NoteOct"), new shared string #16
$ws.Range("E1").Value = "NoteOct"

# New column E: NoteOct = Note & Octave, for data rows 14..129
# (matches the diff's shared-formula groups E14, E15:E78, E79:E129)
$ws.Range("E14").Formula = "=C14&D14"
$ws.Range("E15:E78").Formula = "=C15&D15"
$ws.Range("E79:E129").Formula = "=C79&D79"

# Update the view: select E14:E129 (also clears the stale topLeftCell/selection)
$ws.Range("E14:E129").Select()
